# Fruta / hortaliza, semanal
# Rows 2-13 of the sheet are re-shuffled onto a new weekly schedule:
# the Fecha (D), Volumen (M), Precio minimo/maximo/promedio (N/O/P) and
# Precio $/Kg (S) columns are updated per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44196
$ws.Range("M2").Value = 56
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 3000

$ws.Range("D3").Value = 44902
$ws.Range("M3").Value = 35
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("S3").Value = 2400

$ws.Range("D4").Value = 44175
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 4000

$ws.Range("D5").Value = 44188
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("S5").Value = 3000

$ws.Range("D6").Value = 44907
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 25000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 25000
$ws.Range("S6").Value = 5000

$ws.Range("D7").Value = 44186

$ws.Range("D8").Value = 44931
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("S8").Value = 3600

$ws.Range("D9").Value = 44179
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 4000

$ws.Range("D10").Value = 44181
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 20000
$ws.Range("S10").Value = 4000

$ws.Range("D12").Value = 44914
$ws.Range("N12").Value = 23000
$ws.Range("O12").Value = 23000
$ws.Range("P12").Value = 23000
$ws.Range("S12").Value = 4600

$ws.Range("D13").Value = 44189
$ws.Range("M13").Value = 40
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("S13").Value = 3000
